$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status cells to DONE for newly completed tasks
$ws.Range("B17").Value = "DONE"
$ws.Range("B18").Value = "DONE"
$ws.Range("B19").Value = "DONE"
$ws.Range("B20").Value = "DONE"

# B18:B20 previously had no alignment/border formatting (unlike B17); give
# them the same centered look as the other "DONE" cells
$ws.Range("B18:B20").HorizontalAlignment = -4108
$ws.Range("B18:B20").Borders.LineStyle = 1

# Clear the note text in the merged footer cell (A23:F23)
$ws.Range("A23").Value = ""

# Add a new centered, bold-ish styled block G8:Q9 and merge it
$ws.Range("G8:Q9").HorizontalAlignment = -4108
$ws.Range("G8:Q9").Font.Size = 20
[void]($ws.Range("G8:Q9").MergeCells = $true)

# Select the new block (mirrors what the author had selected when saving)
# and zoom out a bit, matching the recorded view change
[void]$ws.Range("G8:Q9").Select()
$ws.Application.ActiveWindow.Zoom = 70
